$d = $word.ActiveDocument

# 1) Empty paragraph right before the table gets a 12pt (sz/szCs=24 half-pts)
#    paragraph-mark font size.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Font.Size = 12
$p2.Range.Font.SizeBi = 12

# 2) Final (last) paragraph in the body: give it the same 12pt paragraph-mark
#    size, then give its run the same size plus the new text "Não há.".
#    This must happen BEFORE the table is deleted, since deleting the table
#    first leaves the cached Paragraphs/Range objects referring to stale
#    positions in this host.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.Font.Size = 12
$last.Range.Font.SizeBi = 12
$last.Range.InsertAfter("Não há.")

# 3) Remove the whole constraints table.
$d.Tables.Item(1).Delete()

# 4) Drop the now-unused custom table style.
$d.Styles.Item("Table1").Delete()

# 5) Make the (already-portrait) page orientation explicit in pgSz.
$d.PageSetup.Orientation = 0
